$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 891.375
$ws.Range("I19").Value = 586
$ws.Range("J19").Value = 1074.6
$ws.Range("K19").Value = 586
$ws.Range("L19").Value = 1074.6
$ws.Range("M19").Value = -411
$ws.Range("N19").Value = -1424.6

# Row 40
$ws.Range("H40").Value = 4303.467
$ws.Range("I40").Value = 2756.125
$ws.Range("J40").Value = 6071.857
$ws.Range("K40").Value = 2756.125
$ws.Range("L40").Value = 6071.857
$ws.Range("M40").Value = -2581.125
$ws.Range("N40").Value = -6421.857

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

# Row 97
$ws.Range("H97").Value = 2000.5714
$ws.Range("J97").Value = 2000.5714
$ws.Range("L97").Value = 6001.7142
$ws.Range("N97").Value = -6993.7142

# Row 99
$ws.Range("H99").Value = 307
$ws.Range("I99").Value = 364
$ws.Range("K99").Value = 1092
$ws.Range("M99").Value = 406

# Row 132
$ws.Range("H132").Value = 1102.7241
$ws.Range("I132").Value = 1102.7241
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3308.1723
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# Row 135
$ws.Range("H135").Value = 933.3333
$ws.Range("I135").Value = 933.3333
$ws.Range("K135").Value = 8399.9997
$ws.Range("M135").Value = -5864.9997

# Row 137
$ws.Range("H137").Value = 1145.7273
$ws.Range("I137").Value = 550
$ws.Range("J137").Value = 1278.1111
$ws.Range("K137").Value = 1650
$ws.Range("L137").Value = 3834.3333
$ws.Range("M137").Value = 900
$ws.Range("N137").Value = -8934.3333

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4864.4
$ws.Range("J32").Value = 8891
$ws.Range("L32").Value = 8891
$ws.Range("N32").Value = -9465

# Row 61
$ws.Range("H61").Value = 16881.875
$ws.Range("I61").Value = 10010.75
$ws.Range("K61").Value = 10010.75
$ws.Range("M61").Value = -9798.75

# Row 63
$ws.Range("H63").Value = 3510.75
$ws.Range("I63").Value = 3349.3333
$ws.Range("K63").Value = 3349.3333
$ws.Range("M63").Value = -2663.3333

# Row 66
$ws.Range("H66").Value = 3510.75
$ws.Range("I66").Value = 3349.3333
$ws.Range("K66").Value = 16746.6665
$ws.Range("M66").Value = -13314.6665

# Row 74
$ws.Range("H74").Value = 2958.4
$ws.Range("I74").Value = 2958.4
$ws.Range("K74").Value = 2958.4
$ws.Range("M74").Value = -2084.4

# Row 77
$ws.Range("H77").Value = 2958.4
$ws.Range("I77").Value = 2958.4
$ws.Range("K77").Value = 14792
$ws.Range("M77").Value = -10424

# Row 132
$ws.Range("H132").Value = 3625
$ws.Range("I132").Value = 3625
$ws.Range("K132").Value = 10875
$ws.Range("M132").Value = -8345

# Row 136
$ws.Range("H136").Value = 16881.875
$ws.Range("I136").Value = 10010.75
$ws.Range("K136").Value = 30032.25
$ws.Range("M136").Value = -27482.25

$ws = $wb.Worksheets.Item("BSM")
# Row 37
$ws.Range("H37").Value = 1787.5
$ws.Range("I37").Value = 716.6667
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 716.6667
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = -579.6667
$ws.Range("N37").Value = -5274

# Row 99
$ws.Range("H99").Value = 1977.75
$ws.Range("I99").Value = 1574.4
$ws.Range("K99").Value = 1574.4
$ws.Range("M99").Value = -76.40000000000009

# Row 134
$ws.Range("H134").Value = 1040.3334
$ws.Range("I134").Value = 1079.6
$ws.Range("J134").Value = 844
$ws.Range("K134").Value = 3238.8
$ws.Range("L134").Value = 2532
$ws.Range("M134").Value = -703.7999999999997
$ws.Range("N134").Value = -7602

$ws = $wb.Worksheets.Item("CRP")
# Row 120
$ws.Range("H120").Value = 43000
$ws.Range("J120").Value = 43000
$ws.Range("L120").Value = 43000
$ws.Range("N120").Value = -50258

# Row 132
$ws.Range("H132").Value = 2483.1428
$ws.Range("I132").Value = 2397
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 7191
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -4661
$ws.Range("N132").Value = -14060

# Row 134
$ws.Range("H134").Value = 2116.25
$ws.Range("I134").Value = 2140.5
$ws.Range("J134").Value = 1995
$ws.Range("K134").Value = 6421.5
$ws.Range("L134").Value = 5985
$ws.Range("M134").Value = -3886.5
$ws.Range("N134").Value = -11055

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 469.25
$ws.Range("I7").Value = 85
$ws.Range("K7").Value = 255
$ws.Range("M7").Value = -143

# Row 19
$ws.Range("H19").Value = 600
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 600
$ws.Range("K19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -2148

# Row 92
$ws.Range("H92").Value = 1090.1818
$ws.Range("I92").Value = 999.375
$ws.Range("J92").Value = 1332.3334
$ws.Range("K92").Value = 2998.125
$ws.Range("L92").Value = 3997.0002
$ws.Range("M92").Value = -1750.125
$ws.Range("N92").Value = -6493.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3562.5
$ws.Range("I132").Value = 4166.6665
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 12499.9995
$ws.Range("L132").Value = 5250
$ws.Range("M132").Value = -9969.999500000002
$ws.Range("N132").Value = -10310

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 10060.2
$ws.Range("I100").Value = 10060.2
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 10060.2
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

# Row 122
$ws.Range("H122").Value = 5828.2856
$ws.Range("I122").Value = 5133
$ws.Range("K122").Value = 15399
$ws.Range("M122").Value = -12949

# Row 132
$ws.Range("H132").Value = 4999.5
$ws.Range("I132").Value = 4999
$ws.Range("K132").Value = 14997
$ws.Range("M132").Value = -12467

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 74933
$ws.Range("J46").Value = 74933
$ws.Range("L46").Value = 74933
$ws.Range("N46").Value = -75395

# Row 132
$ws.Range("H132").Value = 1260.4
$ws.Range("I132").Value = 1260.4
$ws.Range("K132").Value = 3781.2
$ws.Range("M132").Value = -1251.2

# Row 134
$ws.Range("H134").Value = 74933
$ws.Range("J134").Value = 74933
$ws.Range("L134").Value = 224799
$ws.Range("N134").Value = -229869

# Row 136
$ws.Range("H136").Value = 1065.7142
$ws.Range("I136").Value = 1026.6666
$ws.Range("K136").Value = 3079.9998
$ws.Range("M136").Value = -529.9998000000001

